# Refresh the crypto price/volume snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain text in the sheet (e.g.
# "596.58", "  -0.87%  "), not numbers, so every new value below is written
# with a leading "'" text-qualifier. That keeps Excel from reinterpreting
# digit-and-dot strings (or "%"-suffixed strings) as Number cells; the
# qualifier itself is not stored as part of the cell text.

$ws.Range("D2").Value = "'68.092.57"
$ws.Range("E2").Value = "'  -0.87%  "

$ws.Range("D3").Value = "'3.779.30"
$ws.Range("E3").Value = "'  -2.41%  "

$ws.Range("D5").Value = "'596.58"
$ws.Range("E5").Value = "'  -0.96%  "

$ws.Range("D6").Value = "'168.44"
$ws.Range("E6").Value = "'  -1.87%  "

$ws.Range("D7").Value = "'3.778.66"
$ws.Range("E7").Value = "'  -2.43%  "

$ws.Range("E8").Value = "'  +0.00%  "

$ws.Range("E9").Value = "'  -0.68%  "

$ws.Range("E10").Value = "'  -2.91%  "

$ws.Range("E11").Value = "'  +0.04%  "

$ws.Range("E12").Value = "'  -2.45%  "

$ws.Range("D13").Value = "'0.0000278"
$ws.Range("E13").Value = "'  -2.53%  "

$ws.Range("D14").Value = "'36.63"
$ws.Range("E14").Value = "'  -1.49%  "

$ws.Range("D15").Value = "'4.414.38"
$ws.Range("E15").Value = "'  -2.31%  "

$ws.Range("D16").Value = "'3.778.22"
$ws.Range("E16").Value = "'  -2.35%  "

$ws.Range("D17").Value = "'18.69"
$ws.Range("E17").Value = "'  +1.83%  "

$ws.Range("D18").Value = "'68.011.64"
$ws.Range("E18").Value = "'  -1.04%  "

$ws.Range("D19").Value = "'7.16"
$ws.Range("E19").Value = "'  -3.59%  "

$ws.Range("E20").Value = "'  -0.39%  "

$ws.Range("D21").Value = "'10.54"
$ws.Range("E21").Value = "'  -4.49%  "

$ws.Range("D22").Value = "'467.06"
$ws.Range("E22").Value = "'  -1.24%  "

$ws.Range("E23").Value = "'  -2.45%  "

$ws.Range("E24").Value = "'  -8.47%  "

$ws.Range("D25").Value = "'83.83"
$ws.Range("E25").Value = "'  -0.19%  "

$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "'  -1.66%  "

$ws.Range("D27").Value = "'12.14"
$ws.Range("E27").Value = "'  -1.06%  "

$ws.Range("D28").Value = "'10.36"
$ws.Range("E28").Value = "'  -1.49%  "

$ws.Range("E29").Value = "'  -0.07%  "

$ws.Range("E30").Value = "'  -1.22%  "

$ws.Range("D31").Value = "'3.925.65"
$ws.Range("E31").Value = "'  -2.36%  "

$ws.Range("D32").Value = "'7.57"
$ws.Range("E32").Value = "'  -2.98%  "

$ws.Range("E33").Value = "'  -3.08%  "

$ws.Range("E34").Value = "'  -4.40%  "

$ws.Range("E35").Value = "'  -2.49%  "

$ws.Range("D36").Value = "'3.735.53"
$ws.Range("E36").Value = "'  -2.59%  "

$ws.Range("D37").Value = "'3.75"
$ws.Range("E37").Value = "'  -5.00%  "

$ws.Range("E38").Value = "'  -1.85%  "

$ws.Range("B39").Value = "'Mantle"
$ws.Range("C39").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "'1.01"
$ws.Range("E39").Value = "'  -1.43%  "

$ws.Range("B40").Value = "'Kaspa"
$ws.Range("C40").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.138"
$ws.Range("E40").Value = "'  -1.73%  "

$ws.Range("E42").Value = "'  -0.02%  "

$ws.Range("E43").Value = "'  -2.95%  "

$ws.Range("E44").Value = "'  +0.02%  "

$ws.Range("D45").Value = "'8.67"
$ws.Range("E45").Value = "'  -1.50%  "

$ws.Range("E46").Value = "'  -3.26%  "

$ws.Range("D47").Value = "'405.34"
$ws.Range("E47").Value = "'  -3.46%  "

$ws.Range("D48").Value = "'45.59"
$ws.Range("E48").Value = "'  -2.35%  "

$ws.Range("D49").Value = "'143.64"
$ws.Range("E49").Value = "'  +0.72%  "

$ws.Range("D50").Value = "'0.000273"
$ws.Range("E50").Value = "'  -9.64%  "

$ws.Range("D51").Value = "'39.93"
